$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Retrieved pizzas for first order: subtract used ingredient quantities
$ws.Range("A2").Value = 948
$ws.Range("B2").Value = 854
$ws.Range("C2").Value = 854
$ws.Range("D2").Value = 854
$ws.Range("G2").Value = 948
